$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.062.12"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.877.89"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.62"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4912"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2930"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06617"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").Value = "1.884.60"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.60"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07201"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6677"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.44"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "30.034.45"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007820"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9989"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").Value = "2.122.65"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9963"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.788"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.885"
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.136"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.44"
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.53"
$ws.Range("E26").Value = "  +7.32%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.902"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.388"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.210"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08784"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.999"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05071"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7215"
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.113"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.661"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01842"
$ws.Range("E37").Value = "  +10.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.685"
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("E39").Value = "  -3.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9302"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.772"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4233"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9985"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.24"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.388"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1277"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05707"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.81"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3789"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.305"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.345"
$ws.Range("E51").Value = "  -0.52%  "
